$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1) Refresh the existing "GDP per Capita" series in column E (rows 2-192).
#    The source values are plain digit strings (and a few decimals) that were
#    stored as text in the original workbook, so force text storage with a
#    NumberFormat round-trip (set to Text, write values, then clear the
#    temporary format back off) rather than letting them land as numbers.
# ---------------------------------------------------------------------------
$updates = @{
    2   = "877"
    52  = "1116"
    95  = "1594"
    132 = "2651"
    133 = "2702"
    134 = "2751"
    135 = "2801"
    136 = "2848"
    137 = "2590"
    138 = "3410"
    139 = "3354"
    140 = "3539"
    141 = "3645"
    142 = "3714"
    143 = "4280"
    144 = "4176"
    145 = "4222"
    146 = "4752"
    147 = "5074"
    148 = "5000"
    149 = "4876"
    150 = "4262"
    151 = "4420"
    152 = "3818"
    153 = "3771"
    154 = "3754"
    155 = "3806"
    156 = "3995"
    157 = "4117"
    158 = "4935"
    159 = "5072"
    160 = "5926"
    161 = "6247"
    162 = "7141"
    163 = "7176"
    164 = "7401"
    165 = "7262"
    166 = "7599"
    167 = "7578"
    168 = "7973"
    169 = "7911"
    170 = "7571"
    171 = "6540"
    172 = "6044"
    173 = "5900.50932029287"
    174 = "6456.12314842036"
    175 = "6460.66282413885"
    176 = "6538.17325423484"
    177 = "6760.10846508219"
    178 = "6783.22904950219"
    179 = "6942.458967216"
    180 = "7120.49369152888"
    181 = "7340.49368937669"
    182 = "7621.84061961081"
    183 = "7982.74070250386"
    184 = "8394.05508411191"
    185 = "8669.62464947069"
    186 = "9288.61659605206"
    187 = "9860.80814524952"
    188 = "10404.5617208557"
    189 = "10938.8043833142"
    190 = "11378.2562593997"
    191 = "11651.249751168"
    192 = "11601.0875374148"
}

$updateRange = $ws.Range("E2:E192")
$updateRange.NumberFormat = "@"
foreach ($r in $updates.Keys) {
    $ws.Cells.Item($r, 5).Value = $updates[$r]
}
$updateRange.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Append six new "Data" rows (years 2011-2016) right after the old last
#    row (192 / year 2010), carrying the same Country Code / Country Name /
#    Indicator as every other row.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 193; Year = 2011; Value = "11615" }
    @{ Row = 194; Year = 2012; Value = "11523" }
    @{ Row = 195; Year = 2013; Value = "11488" }
    @{ Row = 196; Year = 2014; Value = "11523" }
    @{ Row = 197; Year = 2015; Value = "11519" }
    @{ Row = 198; Year = 2016; Value = "11529" }
)

$newRange = $ws.Range("E193:E198")
$newRange.NumberFormat = "@"

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 400
    $ws.Cells.Item($r, 2).Value = "Jordan"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $item.Year
    $ws.Cells.Item($r, 5).Value = $item.Value
}

$newRange.ClearFormats()
